# Update 'Total Deaths' (column B) values per pre-processor module refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 2642
    3 = 2399
    4 = 3280
    6 = 603
    8 = 64792
    9 = 4149
    11 = 10260
    12 = 4580
    13 = 210
    14 = 657
    15 = 11644
    17 = 2572
    18 = 24322
    20 = 100
    22 = 13035
    24 = 724
    25 = 408622
    27 = 16548
    31 = 226
    32 = 106
    34 = 24334
    37 = 26659
    39 = 75164
    41 = 147
    43 = 3290
    44 = 287
    45 = 7218
    46 = 675
    47 = 321
    48 = 29365
    52 = 3499
    53 = 18765
    54 = 13531
    55 = 2134
    57 = 12
    58 = 1172
    60 = 3757
    62 = 915
    63 = 105291
    64 = 140
    66 = 4163
    67 = 83605
    69 = 10587
    71 = 7578
    72 = 146
    74 = 303
    76 = 5367
    77 = 27908
    79 = 222408
    80 = 45949
    81 = 72875
    82 = 15566
    84 = 6367
    85 = 121433
    86 = 790
    87 = 10391
    88 = 8925
    89 = 3352
    90 = 2781
    91 = 1840
    93 = 1590
    94 = 1622
    96 = 2154
    97 = 7345
    100 = 3047
    102 = 3956
    103 = 800
    104 = 677
    106 = 1551
    108 = 491
    113 = 217345
    114 = 5850
    116 = 126
    117 = 1510
    118 = 9032
    121 = 3362
    122 = 17443
    127 = 757
    128 = 2053
    129 = 18310
    130 = 6244
    132 = 6653
    133 = 62375
    134 = 17525
    135 = 68105
    137 = 480
    138 = 28380
    139 = 109341
    147 = 6992
    148 = 1111
    149 = 6432
    153 = 11807
    154 = 4269
    157 = 54452
    159 = 78293
    160 = 709
    164 = 10655
    165 = 1610
    169 = 276
    172 = 179
    173 = 10915
    174 = 41191
    175 = 577528
    177 = 46607
    178 = 1596
    179 = 127797
    180 = 2796
    181 = 653
    183 = 2189
    185 = 3283
    186 = 1239
    188 = 1573
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
